{"js": "// Replace the sentence about the subnetting table with the corrected\n// \"subnet mask\" wording, preserving the run's existing formatting\n// (it-IT language tag).\nconst body = context.document.body;\n\nconst oldText = \"La tabella di subnetting per la rete che ho creato \u00e8 la seguente:\";\nconst newText = \"La tabella della subnet mask per la rete che ho creato \u00e8 la seguente:\";\n\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found: \" + oldText);\n}\n\n// Replace the whole matched range's text in place so the run keeps its\n// original formatting (rPr) instead of inheriting formatting from\n// neighbouring runs.\nresults.items[0].insertText(newText, \"Replace\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$oldText = \"La tabella di subnetting per la rete che ho creato \u00e8 la seguente:\"\n$newText = \"La tabella della subnet mask per la rete che ho creato \u00e8 la seguente:\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Forward = $true\n$find.Wrap = 1\n\nif ($find.Execute()) {\n    # $rng now spans exactly the matched text; replacing its .Text keeps\n    # the run's original formatting (rPr) intact.\n    $rng.Text = $newText\n} else {\n    throw \"Target sentence not found: $oldText\"\n}\n"}
